$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" (summary) sheet: insert a new top data row for 2022-Q3 ---
# the sheet keeps its newest-quarter-first history; shift existing rows down
# and insert the new quarter at the top, values below stay as before.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.89
# restore the numeric "rank" style on A2 to match the other data rows (A3..A6)
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# --- 2. Add a new "2022-Q3" worksheet (same layout/style as the other quarterly
#        per-fund holding sheets), positioned right after "总计" / before "2022-Q2" ---
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# overwrite the copied holding figures with the new quarter's numbers
# (kept as text, matching the source sheet's inlineStr cell type)
$wsQ3.Range("D2:G2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "27.03"
$wsQ3.Range("E2").Value = "99.07"
$wsQ3.Range("F2").Value = "3.29"
$wsQ3.Range("G2").Value = "0.8893"
$wsQ3.Range("D2:G2").ClearFormats()

# keep "2021-Q2" as the selected/active sheet, as it was before the edit
$wb.Worksheets.Item("2021-Q2").Activate()
